$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "ISIN"
$ws.Range("B1").Value = "Stock Name"
$ws.Range("C1").Value = "Mutual Fund"
$ws.Range("D1").Value = "Jan_2026"
$ws.Range("E1").Value = "Dec_2025"
$ws.Range("F1").Value = "Nov_2025"
$ws.Range("G1").Value = "MoM"
$ws.Range("H1").Value = "QoQ"
$ws.Range("A2").Value = "INE040A01034"
$ws.Range("B2").Value = "HDFC Bank Limited"
$ws.Range("C2").Value = "quant Business Cycle Fund"
$ws.Range("D2").Value = 9.194146
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 9.194146
$ws.Range("H2").Value = 9.194146
$ws.Range("A3").Value = "INE202B01038"
$ws.Range("B3").Value = "Piramal Finance Ltd"
$ws.Range("C3").Value = "quant Business Cycle Fund"
$ws.Range("D3").Value = 8.702396999999999
$ws.Range("E3").Value = 7.671751
$ws.Range("F3").Value = 6.737916
$ws.Range("G3").Value = 1.030645999999999
$ws.Range("H3").Value = 1.964480999999999
$ws.Range("A4").Value = "INE775A01035"
$ws.Range("B4").Value = "Samvardhana Motherson International Ltd"
$ws.Range("C4").Value = "quant Business Cycle Fund"
$ws.Range("D4").Value = 8.288983
$ws.Range("E4").Value = 8.237802
$ws.Range("F4").Value = 7.443303
$ws.Range("G4").Value = 0.0511809999999997
$ws.Range("H4").Value = 0.8456799999999998
$ws.Range("A5").Value = "INE364U01010"
$ws.Range("B5").Value = "Adani Green Energy Limited"
$ws.Range("C5").Value = "quant Business Cycle Fund"
$ws.Range("D5").Value = 8.12912
$ws.Range("E5").Value = 9.057323
$ws.Range("F5").Value = 4.345319
$ws.Range("G5").Value = -0.9282029999999999
$ws.Range("H5").Value = 3.783801
$ws.Range("A6").Value = "INE406A01037"
$ws.Range("B6").Value = "Aurobindo Pharma Limited"
$ws.Range("C6").Value = "quant Business Cycle Fund"
$ws.Range("D6").Value = 7.148391
$ws.Range("E6").Value = 6.550493
$ws.Range("F6").Value = 6.328911
$ws.Range("G6").Value = 0.5978979999999998
$ws.Range("H6").Value = 0.8194800000000004
$ws.Range("A7").Value = "INE795G01014"
$ws.Range("B7").Value = "HDFC Life Insurance Co Ltd"
$ws.Range("C7").Value = "quant Business Cycle Fund"
$ws.Range("D7").Value = 6.989301
$ws.Range("E7").Value = 6.707038
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0.2822630000000004
$ws.Range("H7").Value = 6.989301
$ws.Range("A8").Value = "INE090A01021"
$ws.Range("B8").Value = "ICICI Bank Limited"
$ws.Range("C8").Value = "quant Business Cycle Fund"
$ws.Range("D8").Value = 5.344555
$ws.Range("E8").Value = 2.055998
$ws.Range("F8").Value = 1.981163
$ws.Range("G8").Value = 3.288557
$ws.Range("H8").Value = 3.363392
$ws.Range("A9").Value = "INE768C01028"
$ws.Range("B9").Value = "Zydus Wellness Ltd"
$ws.Range("C9").Value = "quant Business Cycle Fund"
$ws.Range("D9").Value = 3.578534
$ws.Range("E9").Value = 3.394735
$ws.Range("F9").Value = 2.995026
$ws.Range("G9").Value = 0.183799
$ws.Range("H9").Value = 0.5835079999999997
$ws.Range("A10").Value = "INE423A01024"
$ws.Range("B10").Value = "Adani Enterprises Limited"
$ws.Range("C10").Value = "quant Business Cycle Fund"
$ws.Range("D10").Value = 3.248288
$ws.Range("E10").Value = 3.368582
$ws.Range("F10").Value = 0.895957
$ws.Range("G10").Value = -0.1202939999999999
$ws.Range("H10").Value = 2.352331
$ws.Range("A11").Value = "INE016A01026"
$ws.Range("B11").Value = "Dabur India Limited"
$ws.Range("C11").Value = "quant Business Cycle Fund"
$ws.Range("D11").Value = 2.978884
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 2.978884
$ws.Range("H11").Value = 2.978884
$ws.Range("A12").Value = "INE127D01025"
$ws.Range("B12").Value = "HDFC Asset Management Company Ltd"
$ws.Range("C12").Value = "quant Business Cycle Fund"
$ws.Range("D12").Value = 2.902886
$ws.Range("E12").Value = 2.884563
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0.01832300000000009
$ws.Range("H12").Value = 2.902886
$ws.Range("A13").Value = "INE650L01011"
$ws.Range("B13").Value = "BROOKS Laboratories Limited"
$ws.Range("C13").Value = "quant Business Cycle Fund"
$ws.Range("D13").Value = 1.824076
$ws.Range("E13").Value = 2.221796
$ws.Range("F13").Value = 2.640019
$ws.Range("G13").Value = -0.3977199999999999
$ws.Range("H13").Value = -0.8159430000000001
$ws.Range("A14").Value = "INE331A01037"
$ws.Range("B14").Value = "The Ramco Cements Limited"
$ws.Range("C14").Value = "quant Business Cycle Fund"
$ws.Range("D14").Value = 1.15084
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 1.15084
$ws.Range("H14").Value = 1.15084
$ws.Range("A15").Value = "INE259A01022"
$ws.Range("B15").Value = "Colgate-Palmolive (India) Ltd"
$ws.Range("C15").Value = "quant Business Cycle Fund"
$ws.Range("D15").Value = 0.966873
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0.966873
$ws.Range("H15").Value = 0.966873
$ws.Range("A16").Value = "INE522F01014"
$ws.Range("B16").Value = "Coal India Ltd"
$ws.Range("C16").Value = "quant Business Cycle Fund"
$ws.Range("D16").Value = 0.21769
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0.21769
$ws.Range("H16").Value = 0.21769
$ws.Range("A17").Value = "INE548A01028"
$ws.Range("B17").Value = "HFCL Limited"
$ws.Range("C17").Value = "quant Business Cycle Fund"
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 6.160869
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = -6.160869
$ws.Range("A18").Value = "INE467B01029"
$ws.Range("B18").Value = "Tata Consultancy Services Limited"
$ws.Range("C18").Value = "quant Business Cycle Fund"
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 3.078667
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = -3.078667
$ws.Range("H18").Value = 0
$ws.Range("A19").Value = "INE758E01017"
$ws.Range("B19").Value = "Jio Financial Services Limited"
$ws.Range("C19").Value = "quant Business Cycle Fund"
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 8.585172
$ws.Range("F19").Value = 8.304379000000001
$ws.Range("G19").Value = -8.585172
$ws.Range("H19").Value = -8.304379000000001
$ws.Range("A20").Value = "INE002A01018"
$ws.Range("B20").Value = "Reliance Industries Limited"
$ws.Range("C20").Value = "quant Business Cycle Fund"
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 2.777088
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = -2.777088
$ws.Range("A21").Value = "INE424H01027"
$ws.Range("B21").Value = "SUN TV Network Limited"
$ws.Range("C21").Value = "quant Business Cycle Fund"
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 3.834342
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = -3.834342
$ws.Range("A22").Value = "INE423A20016"
$ws.Range("B22").Value = "Adani Enterprises Limited Rights"
$ws.Range("C22").Value = "quant Business Cycle Fund"
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0.020857
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = -0.020857
$ws.Range("A23").Value = "INE271C01023"
$ws.Range("B23").Value = "DLF Limited"
$ws.Range("C23").Value = "quant Business Cycle Fund"
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 2.841351
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = -2.841351
$ws.Range("H23").Value = 0
$ws.Range("A24").Value = "INE245A01021"
$ws.Range("B24").Value = "Tata Power Company Limited"
$ws.Range("C24").Value = "quant Business Cycle Fund"
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 3.049091
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = -3.049091
$ws.Range("A25").Value = "INE212I01016"
$ws.Range("B25").Value = "S. P. Apparels Limited"
$ws.Range("C25").Value = "quant Business Cycle Fund"
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 5.872413
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = -5.872413
$ws.Range("A26").Value = "INE177F01017"
$ws.Range("B26").Value = "Kovai Medical Center & Hospital Ltd."
$ws.Range("C26").Value = "quant Business Cycle Fund"
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0.125405
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = -0.125405
$ws.Range("A27").Value = "INE14LE01019"
$ws.Range("B27").Value = "Aditya Birla Lifestyle Brands Limited"
$ws.Range("C27").Value = "quant Business Cycle Fund"
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1.153089
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = -1.153089
$ws.Range("A28").Value = "INE094A01015"
$ws.Range("B28").Value = "Hindustan Petroleum Corporation Ltd"
$ws.Range("C28").Value = "quant Business Cycle Fund"
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 2.471008
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = -2.471008
$ws.Range("H28").Value = 0
$ws.Range("A29").Value = "INE939A01011"
$ws.Range("B29").Value = "Strides Pharma Science Ltd"
$ws.Range("C29").Value = "quant Business Cycle Fund"
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 1.67295
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = -1.67295
